$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($ws, $cellref, $value) {
    $c = $ws.Range($cellref)
    $c.NumberFormat = "@"
    $c.Value = $value
    $c.Style = "Normal"
}

Set-TextValue $ws 'D2' '25.970.66'
$ws.Range('E2').Value = '  +0.51%  '
Set-TextValue $ws 'D3' '1.647.96'
$ws.Range('E3').Value = '  +0.87%  '
$ws.Range('E4').Value = '  +0.57%  '
Set-TextValue $ws 'D5' '216.20'
$ws.Range('E5').Value = '  +0.74%  '
Set-TextValue $ws 'D6' '0.5104'
$ws.Range('E6').Value = '  +1.80%  '
Set-TextValue $ws 'D7' '1.006'
$ws.Range('E7').Value = '  +0.49%  '
Set-TextValue $ws 'D8' '0.2581'
$ws.Range('E8').Value = '  +0.66%  '
Set-TextValue $ws 'D9' '0.06430'
$ws.Range('E9').Value = '  +0.71%  '
Set-TextValue $ws 'D10' '19.70'
$ws.Range('E10').Value = '  +0.25%  '
Set-TextValue $ws 'D11' '0.07787'
$ws.Range('E11').Value = '  +1.39%  '
Set-TextValue $ws 'D12' '4.327'
$ws.Range('E12').Value = '  +1.97%  '
Set-TextValue $ws 'D13' '1.650.63'
$ws.Range('E13').Value = '  +0.14%  '
Set-TextValue $ws 'D14' '0.5470'
Set-TextValue $ws 'D15' '0.0₅7892'
$ws.Range('E15').Value = '  -0.32%  '
Set-TextValue $ws 'D16' '64.69'
$ws.Range('E16').Value = '  +1.82%  '
Set-TextValue $ws 'D17' '26.039.60'
$ws.Range('E17').Value = '  +0.75%  '
$ws.Range('E18').Value = '  +0.55%  '
Set-TextValue $ws 'D19' '198.71'
$ws.Range('E19').Value = '  -1.40%  '
Set-TextValue $ws 'D20' '4.478'
$ws.Range('E20').Value = '  +3.41%  '
Set-TextValue $ws 'D21' '10.03'
$ws.Range('E21').Value = '  +1.11%  '
Set-TextValue $ws 'D22' '6.061'
$ws.Range('E22').Value = '  +1.63%  '
Set-TextValue $ws 'D23' '1.009'
$ws.Range('E23').Value = '  +0.69%  '
Set-TextValue $ws 'D24' '1.879'
$ws.Range('E24').Value = '  -2.63%  '
Set-TextValue $ws 'D25' '140.41'
$ws.Range('E25').Value = '  -0.95%  '
Set-TextValue $ws 'D26' '0.1152'
$ws.Range('E26').Value = '  +1.08%  '
$ws.Range('E27').Value = '  +3.03%  '
Set-TextValue $ws 'D28' '15.75'
$ws.Range('E28').Value = '  +0.61%  '
Set-TextValue $ws 'D29' '1.244'
$ws.Range('E29').Value = '  +0.49%  '
$ws.Range('E30').Value = '  +0.57%  '
$ws.Range('E31').Value = '  +0.84%  '
Set-TextValue $ws 'D32' '3.210'
$ws.Range('E32').Value = '  +0.87%  '
Set-TextValue $ws 'D33' '1.549'
$ws.Range('E33').Value = '  +0.58%  '
Set-TextValue $ws 'D34' '2.368'
$ws.Range('E34').Value = '  +0.09%  '
Set-TextValue $ws 'D35' '0.8963'
$ws.Range('E35').Value = '  +0.53%  '
Set-TextValue $ws 'D36' '2.593'
$ws.Range('E36').Value = '  -0.79%  '
Set-TextValue $ws 'D37' '1.136.38'
$ws.Range('E37').Value = '  -3.31%  '
Set-TextValue $ws 'D38' '0.5542'
$ws.Range('E38').Value = '  -0.73%  '
$ws.Range('E39').Value = '  +0.43%  '
$ws.Range('E40').Value = '  +0.72%  '
Set-TextValue $ws 'D41' '2.556'
$ws.Range('E41').Value = '  -0.48%  '
Set-TextValue $ws 'D42' '5.665'
$ws.Range('E42').Value = '  -0.46%  '
Set-TextValue $ws 'D43' '0.8180'
$ws.Range('E43').Value = '  +1.50%  '
$ws.Range('B44').Value = 'BabyDogeCoin'
$ws.Range('C44').Value = 'https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge'
Set-TextValue $ws 'D44' '0.0₈125'
$ws.Range('E44').Value = '  +9.26%  '
$ws.Range('B45').Value = 'Quant'
$ws.Range('C45').Value = 'https://coinranking.com/coin/bauj_21eYVwso+quant-qnt'
Set-TextValue $ws 'D45' '100.06'
$ws.Range('E45').Value = '  +0.58%  '
Set-TextValue $ws 'D46' '1.786.89'
$ws.Range('E46').Value = '  +0.96%  '
Set-TextValue $ws 'D47' '0.4535'
$ws.Range('E47').Value = '  +0.50%  '
$ws.Range('B48').Value = 'Aave'
$ws.Range('C48').Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
Set-TextValue $ws 'D48' '55.27'
$ws.Range('E48').Value = '  +0.95%  '
$ws.Range('B49').Value = 'Frax'
$ws.Range('C49').Value = 'https://coinranking.com/coin/KfWtaeV1W+frax-frax'
Set-TextValue $ws 'D49' '1.007'
$ws.Range('E49').Value = '  +0.44%  '
Set-TextValue $ws 'D50' '0.05093'
$ws.Range('E50').Value = '  +0.30%  '
$ws.Range('E51').Value = '  +0.67%  '
